$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = "Flowise"
$ws.Range("A16").Select()
